$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4-8: PunishGold/Food-cycle "count" column C changes from 1 to 0 ---
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0

# Row 4: RewardMental (M4) column no longer used for this row - clear it
$ws.Range("M4").ClearContents()

# --- New string-introducing writes, issued in the order that yields the
#     target shared-string table layout (wolfnest, 狼穴, 赌博, gamble,
#     RewardItem, RewardDrop, 鱼塘, fishpool) ---

# Row 9 (id 42000006) - "wolfnest" (狼穴)
$ws.Range("D9").Value = "wolfnest"
$ws.Range("B9").Value = "狼穴"

# Row 10 (id 42000007) - "gamble" (赌博)
$ws.Range("B10").Value = "赌博"
$ws.Range("E10").Value = "gamble"

# Rename reward-item table columns (header row 1)
$ws.Range("L1").Value = "RewardItem"
$ws.Range("M1").Value = "RewardDrop"

# Row 11 (id 42000008) - "fishpool" (鱼塘)
$ws.Range("B11").Value = "鱼塘"
$ws.Range("D11").Value = "fishpool"

# --- Fill in the remaining cells for the new rows ---

# Row 9
$ws.Range("A9").Value = 42000006
$ws.Range("C9").Value = 0
$ws.Range("E9").Value = "wolfnest"
$ws.Range("F9").Value = 43000002
$ws.Range("G9").Value = 100
$ws.Range("K9").Value = 50
$ws.Range("P9").Value = 100

# Row 10
$ws.Range("A10").Value = 42000007
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = "wolfnest"
$ws.Range("G10").Value = 250
$ws.Range("N10").Value = 250
$ws.Range("Q10").Value = 50

# Row 11
$ws.Range("A11").Value = 42000008
$ws.Range("C11").Value = 0
$ws.Range("E11").Value = "fishpool"
$ws.Range("H11").Value = 100
$ws.Range("M11").Value = 23000001

# --- Expand the table / autofilter range to cover the new rows ---
$tbl = $ws.ListObjects.Item(1)
[void]$tbl.Resize($ws.Range("A1:Q11"))

# --- Match the recorded selection state ---
[void]$ws.Range("H11").Select()
